#
# feat: Selectors- element specific
#
# Adds a new slide ("Class and ID Selector") after the existing "ID Selector"
# slide (slide 9 / sldId 264), using the same "Title and Content" layout.
#

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helpers for building a TextRange out of multiple paragraphs / runs while
# keeping each appended chunk correctly scoped to its own run (so per-run
# character formatting - e.g. color - only applies to that chunk).
# ---------------------------------------------------------------------------

function Set-FirstParagraph($range, $text) {
    $range.Text = $text
    $newLen = $range.Length
    $start = $newLen - $text.Length + 1
    return $range.Characters($start, $text.Length)
}

function Add-Paragraph($range, $text) {
    $range.InsertAfter("`r" + $text) | Out-Null
    $newLen = $range.Length
    $start = $newLen - $text.Length + 1
    return $range.Characters($start, $text.Length)
}

function Add-Run($range, $text) {
    $range.InsertAfter($text) | Out-Null
    $newLen = $range.Length
    $start = $newLen - $text.Length + 1
    return $range.Characters($start, $text.Length)
}

function Set-LastParagraphIndent($range, $level) {
    $count = $range.Paragraphs().Count
    $para = $range.Paragraphs($count, 1)
    $para.IndentLevel = $level
}

# ---------------------------------------------------------------------------
# Add the new slide, reusing the "Title and Content" layout already used by
# the other content slides in the deck (same layout as slide 9).
# ---------------------------------------------------------------------------

$refLayout = $p.Slides.Item(9).CustomLayout
$slide = $p.Slides.AddSlide($p.Slides.Count + 1, $refLayout)

$title = $slide.Shapes.Item(1)
$content = $slide.Shapes.Item(2)

# Title placement matches the other title placeholders on this deck.
$title.Left = 53.3334
$title.Top = 48.0
$title.Width = 676.903
$title.Height = 51.0

$title.TextFrame.TextRange.Text = "Class and ID Selector"

$tr = $content.TextFrame.TextRange

Set-FirstParagraph $tr "One of the strengths of CSS that we have talked about is its ability to cascade, to combine CSS from multiple sources, enable inheritance of styles from parent to child theme." | Out-Null

Add-Paragraph $tr "IDs can only be used once on a page, their scope is more specific and they can override class based styles." | Out-Null

Add-Paragraph $tr "Example:" | Out-Null

# --- Example paragraph 1: HTML markup (red), level 4 ---
$run = Add-Paragraph $tr "HTML: <div id=`u201ccontent`u201d class=`u201cmain`u201d> "
$run.Font.Color.RGB = 255
Set-LastParagraphIndent $tr 4

$run = Add-Run $tr "Lopem"
$run.Font.Color.RGB = 255

$run = Add-Run $tr " ipsum </div>"
$run.Font.Color.RGB = 255

# --- Example paragraph 2: CSS id rule (red), level 4 ---
$run = Add-Paragraph $tr "CSS:  #content { "
$run.Font.Color.RGB = 255
Set-LastParagraphIndent $tr 4

$run = Add-Run $tr "color:green"
$run.Font.Color.RGB = 255

$run = Add-Run $tr " }"
$run.Font.Color.RGB = 255

# --- Example paragraph 3: CSS class rule (red), level 4 ---
$run = Add-Paragraph $tr "CSS: .main { color: red}"
$run.Font.Color.RGB = 255
Set-LastParagraphIndent $tr 4
